# Append two new scraped postings to the "ランサーズ" sheet, shifting the
# existing rows down, refresh the retrieval timestamp on every row, widen
# column H, and keep the hyperlinks on column F pointing at the right URLs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2025-12-18 01:52:46"

# Row data (rows 2..16) after the edit: Title, Category, Price, Deadline, Url, Score, Skills
$rows = @(
    @("【急募】生成AI×業務効率化の実装を支援するエンジニア募集", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5444662", 385, "🔥AI,Ai ◆効率化"),
    @("産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5450864", 383, "🔥AI,Ai ◆開発"),
    @("初回 既存システムのRuby、Ruby on Railsバージョンアップ及び追加改修", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5456434", 318, "🔥AI,Ai"),
    @("【Zapier設定のみ!作業時間~2時間】スプレッドシート・Gドライブ自動化構築(設計済)", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5456066", 255, "🔥API ◆自動化"),
    @("【スマホアプリ開発】 音声データ推定アプリの依頼", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5456360", 175, "★スマホアプリ ◆開発 ◇アプリ"),
    @("【急募】新規システム開発に伴う要件定義依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5455415", 113, "◆開発,システム開発"),
    @("【急募】YouTube自動化チャンネルの台本生成ワークフロー構築パートナー", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5456199", 88, "◆自動化"),
    @("初回 高度な商用SaaSの新規開発と保守業務", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5455862", 75, "◆開発"),
    @("【急募】MVNO会員向けマイページ新規開発エンジニア募集", "システム開発", "1,000,000 円 ~ 3,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5455513", 75, "◆開発"),
    @("ヤフーオークションで複数の欲しい商品を一括検索するツールの作成", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5455714", 68, "◆ツール"),
    @("【小規模・短納期・急募】アプリMatrixifyを用いたデータ移行検証・マッピング担当募集", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5455675", 38, "◇アプリ"),
    @("wordpressレンダリングを妨げるリソースの除外", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5016989", 33, "○WordPress"),
    @("【改善提案募集】事業管理スプレッドシートの見直し・改善提案をお願いします。", "システム開発", "1,000 ~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5455422", 30, "◇管理"),
    @("【急募】LINEシステム構築・保守運用のプロフェッショナルを求む!", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5456063", 33, $null),
    @("グループ毎の日報をリアルタイムでまとめたい", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5456195", 18, $null)
)

# Remove the old hyperlinks first; they will be rebuilt below once the URLs
# are in their new positions.
$ws.Hyperlinks.Delete()

$r = 2
foreach ($item in $rows) {
    $ws.Range("A$r").Value = $timestamp
    $ws.Range("B$r").Value = $item[0]
    $ws.Range("C$r").Value = $item[1]
    $ws.Range("D$r").Value = $item[2]
    $ws.Range("E$r").Value = $item[3]
    $ws.Range("F$r").Value = $item[4]
    $ws.Range("G$r").Value = $item[5]
    if ($item[6]) {
        $ws.Range("H$r").Value = $item[6]
    } else {
        $ws.Range("H$r").ClearContents()
    }
    $r = $r + 1
}

# Rebuild hyperlinks on column F for every data row, preserving the
# "Hyperlink" cell style (blue/underline, same as before the edit).
for ($i = 2; $i -le 16; $i++) {
    $cell = $ws.Range("F$i")
    $url = $cell.Value2
    $ws.Hyperlinks.Add($cell, $url)
    $cell.Style = "Hyperlink"
}

# Column H needs to grow from width 13 to width 18 to fit the new skill text.
$ws.Columns.Item(8).ColumnWidth = 18 - 0.83
